# "udalost je radek s hodnotou v kteremkoliv sloupci"
# (an event is a row with a value in any of its columns)
#
# Clear the leftover/duplicate cell values that do not belong to their
# row's event anymore - i.e. values that were accidentally left behind
# while the rest of that particular row's event data had already been
# cleared out.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B5").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("D5").ClearContents()

$ws.Range("B15").ClearContents()
$ws.Range("C15").ClearContents()
$ws.Range("E15").ClearContents()

$ws.Range("C19").ClearContents()
$ws.Range("D19").ClearContents()
$ws.Range("E19").ClearContents()

$ws.Range("B22").ClearContents()
$ws.Range("D22").ClearContents()
$ws.Range("E22").ClearContents()
